$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R1 Regression")

# Row 4: M4/N4 change from old address to new address
$ws.Range("M4").Value = "UNIT 35, 146-152 PARRAMATTA RD, HOMEBUSH NSW 2140"
$ws.Range("N4").Value = "UNIT 35, 146-152 PARRAMATTA RD, HOMEBUSH NSW 2140"

# Row 5: add M5/N5 with the new address
$ws.Range("M5").Value = "UNIT 35, 146-152 PARRAMATTA RD, HOMEBUSH NSW 2140"
$ws.Range("N5").Value = "UNIT 35, 146-152 PARRAMATTA RD, HOMEBUSH NSW 2140"

# Update the sheet view / selection: drop the frozen/scrolled topLeftCell and change selection to A3
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A3").Select()
